# C5-PowerPoint.pptx edit script
# 1. Change the table style id on the table in slide 6.
# 2. Re-colour the presentation theme from "Integral" to "Office Theme"
#    (the table style + theme swap correspond to a Design-tab theme
#    change made in the original commit).

$p = $ppt.ActivePresentation

# --- 1. Table style change (slide 6, shape 2 is the table) ---
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{47EDC479-AF55-4293-B58C-E3288E792B8C}")

# --- 2. Theme colour scheme change: Integral -> Office ---
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Item(1).RGB  = 0          # dk1      000000
$colors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388    # dk2      44546A
$colors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2  ED7D31
$colors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB  = 49407      # accent4  FFC000
$colors.Item(9).RGB  = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456    # accent6  70AD47
$colors.Item(11).RGB = 12673797   # hlink    0563C1
$colors.Item(12).RGB = 7491477    # folHlink 954F72
